$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.710559844970703
$ws.Range("B1").Value = 3.097336530685425
$ws.Range("C1").Value = 4.841117858886719
$ws.Range("D1").Value = 1.109140515327454
$ws.Range("E1").Value = 0.6284382343292236
